# TestData.xlsx modifications:
#  - add "LinkedList" worksheet (with Code/expectedOutcome sample rows)
#  - add an "output"/validation-message column to the "python DS" and "Login" sheets
#  - adjust column widths / selections to match the saved workbook state
#  - make "Login" the active tab

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Add the new "LinkedList" sheet at the end of the tab strip
# ------------------------------------------------------------------
$lastIndex  = $wb.Worksheets.Count
$lastSheet  = $wb.Worksheets.Item($lastIndex)
$wsLinked   = $wb.Worksheets.Add($null, $lastSheet)
$wsLinked.Name = "LinkedList"

# ------------------------------------------------------------------
# 2. "python DS" sheet - add an "expected Outcome" column
#    (write cells in B1,B2,B3 order so shared-string ids line up)
# ------------------------------------------------------------------
$wsPy = $wb.Worksheets.Item("python DS")
$wsPy.Range("B1").Value = "expected Outcome"
$wsPy.Range("B2").Value = "popuperror message containing Syntaxerror"
$wsPy.Range("B3").Value = "the user is able to see the output in the console"

# ------------------------------------------------------------------
# 3. "LinkedList" sheet content
#    (write B1,B2,B3 right after A-column so new strings keep this order)
# ------------------------------------------------------------------
$wsLinked.Range("A1").Value = "Code"
$wsLinked.Range("B1").Value = "expectedOutcome"
$wsLinked.Range("A2").Value = 'print(' + [char]34 + 'hello)abc'
$wsLinked.Range("B2").Value = "SyntaxError: bad input on line 1"
$wsLinked.Range("A3").Value = 'print(' + [char]34 + 'hello' + [char]34 + ')'
$wsLinked.Range("B3").Value = "hello"

# ------------------------------------------------------------------
# 4. "Login" sheet - add the "output" column (C2..C8 first, C1 last,
#    matching the original authoring order of the workbook)
# ------------------------------------------------------------------
$wsLogin = $wb.Worksheets.Item("Login")
$wsLogin.Range("C2").Value = "Please fill out this field."
$wsLogin.Range("C3").Value = "Please fill out this field."
$wsLogin.Range("C4").Value = "Please fill out this field."
$wsLogin.Range("C5").Value = "Invalid Username and Password"
$wsLogin.Range("C6").Value = "Invalid Username and Password"
$wsLogin.Range("C7").Value = "Invalid Username and Password"
$wsLogin.Range("C8").Value = "You are logged in"
$wsLogin.Range("C1").Value = "          output"

# ------------------------------------------------------------------
# 5. Column widths
# ------------------------------------------------------------------
$wsPy.Columns.Item(2).ColumnWidth     = 43.3
$wsLinked.Columns.Item(1).ColumnWidth = 18.9
$wsLinked.Columns.Item(2).ColumnWidth = 32.18
$wsLogin.Columns.Item(3).ColumnWidth  = 28.42

# ------------------------------------------------------------------
# 6. Page setup
# ------------------------------------------------------------------
$wsLinked.PageSetup.Orientation = 1

# ------------------------------------------------------------------
# 7. Selections on each sheet
# ------------------------------------------------------------------
$null = $wsPy.Range("D4").Select()
$null = $wsLinked.Range("B5").Select()
$null = $wsLogin.Range("C1").Select()

# ------------------------------------------------------------------
# 8. Make "Login" the active sheet/tab (do this last so it "sticks")
# ------------------------------------------------------------------
$null = $wsLogin.Activate()
$null = $wsLogin.Range("C1").Select()
